# The reviewer ("MArcel") pointed out that the extra parenthetical note
# ("atenção para a diferença de conteúdo entre projeto e pré-projeto")
# that was tacked onto the "REVISÃO BIBLIOGRÁFICA" checklist item is
# redundant/confusing between the PreProjeto and Projeto forms, so it is
# removed from both occurrences in the document (one is a single run,
# the other is split across two runs - Find/Replace across the whole
# story normalizes both back down to a single run with the short text).

$d = $word.ActiveDocument

$oldText = "REVISÃO BIBLIOGRÁFICA (atenção para a diferença de conteúdo entre projeto e pré-projeto)"
$newText = "REVISÃO BIBLIOGRÁFICA"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdFindContinue=1, wdReplaceAll=2 - replaces every match in the document
# (both the single-run and the two-run occurrence collapse to the same
# final text).
$find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, `
              $false, $newText, 2) | Out-Null

# Safety net in case any fragment (e.g. the trailing " (atenção ...)" run
# that follows an already-short "REVISÃO BIBLIOGRÁFICA" run) survived the
# pass above as its own separate match.
$leftover = " (atenção para a diferença de conteúdo entre projeto e pré-projeto)"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute($leftover, $true, $false, $false, $false, $false, $true, 1, `
               $false, "", 2) | Out-Null
